$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New s_val data (regenerated to filter save games)
$data = @{
    2 = @(1.459612070389937, 1.667794583268128, 337.1190423067083, 8.660232485948974, 1, 348.9066814463153)
    3 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
    4 = @(0.01514828764759746, 0.002777888934908601, 0.8054896365839992, 8.660232485948974, 0, 9.48364829911548)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    6 = @(0.127881588408715, 0.002777888934908601, 0.1575252929769615, 0.496779210170732, 0, 0.7849639804913171)
    7 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 8.660232485948974, 0, 39.7764191927396)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
